# Applies the "sent a copy to Nidhi" edit:
#  - corrects three mis-attributed Importer_Name values (rows 14, 32; row 34 unchanged text)
#  - highlights the corrected rows (14, 32, 34) with a solid yellow fill
#  - refreshes the USD-converted columns (T/U/V) that depend on the corrected
#    exchange-rate lookups for the affected rows and the rows that share the
#    same BE_Date / currency grouping

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Correct the Importer_Name (column P) text for the two mis-matched rows.
# ---------------------------------------------------------------------------
$ws.Range("P14").Value = "zuari cement limited"
$ws.Range("P32").Value = "continental india private limited"

# ---------------------------------------------------------------------------
# 2. Highlight the corrected rows in yellow (solid fill, A:V) so they are easy
#    to spot for reviewers. Build the fill on a single cell first, then
#    propagate it with a format-only paste so every row lands on the same
#    style index instead of minting a fresh one per row.
# ---------------------------------------------------------------------------
$source = $ws.Range("A14")
$source.Interior.Color = 65535
$source.Interior.PatternColor = 65535

$source.Copy() | Out-Null
$ws.Range("B14:V14").PasteSpecial(-4122) | Out-Null
$ws.Range("A32:V32").PasteSpecial(-4122) | Out-Null
$ws.Range("A34:V34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Recalculated USD columns (Unit_Price_USD / TOTAL_ASS_VALUE_USD /
#    Invoice_Unit_Price_FC_USD) for the rows whose conversion rates moved.
# ---------------------------------------------------------------------------
$ws.Range("T12").Value = 2.2005
$ws.Range("U12").Value = 2860.6286
$ws.Range("V12").Value = 3.4818

$ws.Range("T15").Value = 6.3817
$ws.Range("U15").Value = 6381471.8322

$ws.Range("T16").Value = 6.3817
$ws.Range("U16").Value = 12762943.3922

$ws.Range("T17").Value = 6.3817
$ws.Range("U17").Value = 15953678.8999

$ws.Range("T18").Value = 6.3817
$ws.Range("U18").Value = 1749448.5212

$ws.Range("T19").Value = 6.3817
$ws.Range("U19").Value = 3190735.78

$ws.Range("T31").Value = 2.1819
$ws.Range("U31").Value = 3491.0509
$ws.Range("V31").Value = 3.4818

$ws.Range("T34").Value = 84.8843
$ws.Range("U34").Value = 848818.9044999999
$ws.Range("V34").Value = 0.8065

$ws.Range("T35").Value = 82.98909999999999
$ws.Range("U35").Value = 331956.3307
$ws.Range("V35").Value = 0.8065

$ws.Range("T36").Value = 82.98909999999999
$ws.Range("U36").Value = 331956.3307
$ws.Range("V36").Value = 0.8065

$ws.Range("T37").Value = 671.5915
$ws.Range("U37").Value = 134317.2859
$ws.Range("V37").Value = 8.2247

$ws.Range("T38").Value = 79873.7458
$ws.Range("U38").Value = 1916969.9
$ws.Range("V38").Value = 984.2729

$ws.Range("T43").Value = 262171.2558
$ws.Range("U43").Value = 1310856.2557
$ws.Range("V43").Value = 2655.7896
